$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "singap"

$ws.Range("C10").Select()

$wb.Save()
